# Auto-generated edit script applying numeric updates to the
# "currentAveragePrice*" / "LevePrice*" / "LeveProfit*" columns (H:N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets,
# per the scheduled pricing-data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2926.7778
$ws.Range("J40").Value = 3127.2856
$ws.Range("L40").Value = 3127.2856
$ws.Range("N40").Value = -3477.2856
$ws.Range("H70").Value = 11781.363
$ws.Range("I70").Value = 3124
$ws.Range("J70").Value = 16728.428
$ws.Range("K70").Value = 9372
$ws.Range("L70").Value = 50185.284
$ws.Range("M70").Value = -9102
$ws.Range("N70").Value = -50725.284
$ws.Range("H73").Value = 11781.363
$ws.Range("I73").Value = 3124
$ws.Range("J73").Value = 16728.428
$ws.Range("K73").Value = 9372
$ws.Range("L73").Value = 50185.284
$ws.Range("M73").Value = -8436
$ws.Range("N73").Value = -52057.284
$ws.Range("H76").Value = 4439.143
$ws.Range("I76").Value = 4284.6665
$ws.Range("K76").Value = 4284.6665
$ws.Range("M76").Value = -3969.6665
$ws.Range("H79").Value = 4439.143
$ws.Range("I79").Value = 4284.6665
$ws.Range("K79").Value = 4284.6665
$ws.Range("M79").Value = -3192.6665
$ws.Range("H98").Value = 10216
$ws.Range("J98").Value = 17833.166
$ws.Range("L98").Value = 17833.166
$ws.Range("N98").Value = -20829.166
$ws.Range("H112").Value = 4342.1514
$ws.Range("J112").Value = 4496.107
$ws.Range("L112").Value = 13488.321
$ws.Range("N112").Value = -15704.321
$ws.Range("H116").Value = 3999
$ws.Range("I116").Value = 3999
$ws.Range("K116").Value = 3999
$ws.Range("M116").Value = -557
$ws.Range("H122").Value = 10216
$ws.Range("J122").Value = 17833.166
$ws.Range("L122").Value = 53499.49800000001
$ws.Range("N122").Value = -58399.49800000001
$ws.Range("H127").Value = 1206.2
$ws.Range("I127").Value = 1132.875
$ws.Range("J127").Value = 1499.5
$ws.Range("K127").Value = 3398.625
$ws.Range("L127").Value = 4498.5
$ws.Range("M127").Value = 1561.375
$ws.Range("N127").Value = -14418.5
$ws.Range("H134").Value = 27228.916
$ws.Range("J134").Value = 27228.916
$ws.Range("L134").Value = 27228.916
$ws.Range("N134").Value = -37368.916
$ws.Range("H138").Value = 4426.1025
$ws.Range("J138").Value = 4644
$ws.Range("L138").Value = 13932
$ws.Range("N138").Value = -24212
$ws.Range("H141").Value = 9998
$ws.Range("I141").Value = 6997
$ws.Range("K141").Value = 20991
$ws.Range("M141").Value = -15811

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1943.2344
$ws.Range("I32").Value = 1469.0161
$ws.Range("K32").Value = 1469.0161
$ws.Range("M32").Value = -1182.0161
$ws.Range("H63").Value = 2674.9412
$ws.Range("I63").Value = 2808.4614
$ws.Range("J63").Value = 2241
$ws.Range("K63").Value = 2808.4614
$ws.Range("L63").Value = 2241
$ws.Range("M63").Value = -2122.4614
$ws.Range("N63").Value = -3613
$ws.Range("H66").Value = 2674.9412
$ws.Range("I66").Value = 2808.4614
$ws.Range("J66").Value = 2241
$ws.Range("K66").Value = 14042.307
$ws.Range("L66").Value = 11205
$ws.Range("M66").Value = -10610.307
$ws.Range("N66").Value = -18069
$ws.Range("H97").Value = 2152.0715
$ws.Range("I97").Value = 1476.4445
$ws.Range("K97").Value = 1476.4445
$ws.Range("M97").Value = -980.4445000000001
$ws.Range("H125").Value = 138847
$ws.Range("J125").Value = 138847
$ws.Range("L125").Value = 138847
$ws.Range("N125").Value = -148687

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 37927.43
$ws.Range("I94").Value = 13098.4
$ws.Range("K94").Value = 13098.4
$ws.Range("M94").Value = -12647.4
$ws.Range("H99").Value = 1853.25
$ws.Range("I99").Value = 808.3333
$ws.Range("K99").Value = 808.3333
$ws.Range("M99").Value = 689.6667
$ws.Range("H133").Value = 96990
$ws.Range("J133").Value = 96990
$ws.Range("L133").Value = 96990
$ws.Range("N133").Value = -107110
$ws.Range("H134").Value = 6572.521
$ws.Range("I134").Value = 6485.844
$ws.Range("K134").Value = 19457.532
$ws.Range("M134").Value = -16922.532

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 19566.25
$ws.Range("I62").Value = 7997.8
$ws.Range("J62").Value = 38847
$ws.Range("K62").Value = 7997.8
$ws.Range("L62").Value = 38847
$ws.Range("M62").Value = -7373.8
$ws.Range("N62").Value = -40095
$ws.Range("H65").Value = 19566.25
$ws.Range("I65").Value = 7997.8
$ws.Range("J65").Value = 38847
$ws.Range("K65").Value = 39989
$ws.Range("L65").Value = 194235
$ws.Range("M65").Value = -36869
$ws.Range("N65").Value = -200475
$ws.Range("H105").Value = 3312.6667
$ws.Range("I105").Value = 2638.1667
$ws.Range("J105").Value = 4661.6665
$ws.Range("K105").Value = 2638.1667
$ws.Range("L105").Value = 4661.6665
$ws.Range("M105").Value = -891.1667000000002
$ws.Range("N105").Value = -8155.6665
$ws.Range("H132").Value = 3398
$ws.Range("I132").Value = 2396.25
$ws.Range("K132").Value = 7188.75
$ws.Range("M132").Value = -4658.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 185.125
$ws.Range("J2").Value = 195
$ws.Range("L2").Value = 1170
$ws.Range("N2").Value = -1396
$ws.Range("H38").Value = 190.23077
$ws.Range("I38").Value = 102.4
$ws.Range("K38").Value = 307.2
$ws.Range("M38").Value = 39.79999999999995
$ws.Range("H122").Value = 1499
$ws.Range("J122").Value = 2226.182
$ws.Range("L122").Value = 20035.638
$ws.Range("N122").Value = -24935.638

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2104400.2
$ws.Range("I14").Value = 3503333.8
$ws.Range("K14").Value = 3503333.8
$ws.Range("M14").Value = -3503165.8
$ws.Range("H102").Value = 11477.186
$ws.Range("I102").Value = 13329.143
$ws.Range("K102").Value = 13329.143
$ws.Range("M102").Value = -11707.143
$ws.Range("H126").Value = 8926
$ws.Range("I126").Value = 9989.333000000001
$ws.Range("J126").Value = 7650
$ws.Range("K126").Value = 29967.999
$ws.Range("L126").Value = 22950
$ws.Range("M126").Value = -27497.999
$ws.Range("N126").Value = -27890
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2485.25
$ws.Range("I40").Value = 2111.625
$ws.Range("K40").Value = 2111.625
$ws.Range("M40").Value = -1975.625
$ws.Range("H122").Value = 6215.25
$ws.Range("I122").Value = 6517.8
$ws.Range("K122").Value = 19553.4
$ws.Range("M122").Value = -17103.4
$ws.Range("H132").Value = 2848.6553
$ws.Range("I132").Value = 2272.8572
$ws.Range("K132").Value = 6818.571599999999
$ws.Range("M132").Value = -4288.571599999999
$ws.Range("H136").Value = 2094.6206
$ws.Range("I136").Value = 1933
$ws.Range("J136").Value = 2401.7
$ws.Range("K136").Value = 5799
$ws.Range("L136").Value = 7205.099999999999
$ws.Range("M136").Value = -3249
$ws.Range("N136").Value = -12305.1

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8984.5
$ws.Range("J62").Value = 8999
$ws.Range("L62").Value = 8999
$ws.Range("N62").Value = -10247
$ws.Range("H65").Value = 8984.5
$ws.Range("J65").Value = 8999
$ws.Range("L65").Value = 44995
$ws.Range("N65").Value = -51235
